$d = $word.ActiveDocument

# --- 1. Main body: "TERE" (unique, bold run after "A ") -> "QWER" ---
$bodyRng = $d.Content
$bodyRng.Find.Execute("TERE", $true, $false, $false, $false, $false, $true, 1, $false, "QWER", 1) | Out-Null

# --- 2. Header: several "Tre"-family runs, each replaced with its own text ---
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$hdrRng = $hdr.Range

$hdrReplacements = @(
    @{ Find = "TRE";  Replace = "QWER" },
    @{ Find = "TERE"; Replace = "QWER" },
    @{ Find = "Tre";  Replace = "Qwer" },
    @{ Find = "Tre";  Replace = "Qwer" },
    @{ Find = "Tre";  Replace = "Qewr" },
    @{ Find = "Tre";  Replace = "Qewr" },
    @{ Find = "Tre";  Replace = "Qwer" },
    @{ Find = "tre";  Replace = "qwer" },
    @{ Find = "tre";  Replace = "qwer" },
    @{ Find = "tre";  Replace = "qwer" }
)

foreach ($item in $hdrReplacements) {
    $hdrRng.Find.Execute($item.Find, $true, $false, $false, $false, $false, $true, 1, $false, $item.Replace, 1) | Out-Null
}
